# 0.02.10 - Add check groups for steel members
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row (row 1): append "*" to all existing parameter headers
#    (B1:M1) and add new headers for the check-group columns N1:S1.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "x1*"
$ws.Range("C1").Value = "y1*"
$ws.Range("D1").Value = "z1*"
$ws.Range("E1").Value = "pin1*"
$ws.Range("F1").Value = "sectionType*"
$ws.Range("G1").Value = "sectionName*"
$ws.Range("H1").Value = "betaAngle*"
$ws.Range("I1").Value = "isDivided*"
$ws.Range("J1").Value = "pin2*"
$ws.Range("K1").Value = "x2*"
$ws.Range("L1").Value = "y2*"
$ws.Range("M1").Value = "z2*"

$ws.Range("N1").Value = "steel*"
$ws.Range("O1").Value = "Ry*"
$ws.Range("P1").Value = "muXZ*"
$ws.Range("Q1").Value = "muXY*"
$ws.Range("R1").Value = "gammaC*"
$ws.Range("S1").Value = "FC*"

# Copy the header formatting (style s=2, same as A1) onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("N1:S1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 2) Group separator row 2 ("Columns"): extend formatting into N2:O2
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("N2:O2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 3) Relabel a couple of member names
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "C4"
$ws.Range("A8").Value = "Балка 1"

# ---------------------------------------------------------------------
# 4) Fill in the new "check group" columns (N:S) for every data row.
#    N = steel grade (text), O = Ry (number), P = muXZ, Q = muXY,
#    R = gammaC, S = FC
# ---------------------------------------------------------------------
function Set-CheckGroup($row, $Ry, $muXZ, $muXY, $gammaC, $FC) {
    $ws.Range("N$row").Value = "C255"
    $ws.Range("O$row").Value = $Ry
    $ws.Range("P$row").Value = $muXZ
    $ws.Range("Q$row").Value = $muXY
    $ws.Range("R$row").Value = $gammaC
    $ws.Range("S$row").Value = $FC
}

Set-CheckGroup 3  240 2.5 3 1    120
Set-CheckGroup 4  240 3.5 3 2    125
Set-CheckGroup 5  240 4.5 3 3    130
Set-CheckGroup 6  240 5.5 3 4    135
Set-CheckGroup 8  240 1   0.5 0.95 400
Set-CheckGroup 9  240 1   1   0.95 400
Set-CheckGroup 11 240 1   0.9 1.05 200
Set-CheckGroup 12 240 0.89 1  1.05 200

# Apply correct cell formatting to the new columns:
#  - N, P, Q, R, S use the regular numeric/text style (same as column B, s=1)
#  - O uses the "isDivided"-like integer style (same as column I, s=6)
$ws.Range("B3").Copy()
$ws.Range("N3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P3:S3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P4:S4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P5:S5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P6:S6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P8:S8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P9:S9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P11:S11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P12:S12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("I3").Copy()
$ws.Range("O3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 5) Column widths: column F (6) grows slightly, and the new column O
#    (15) gets an explicit custom width.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(15).ColumnWidth = 8.6666666666666667

# ---------------------------------------------------------------------
# 6) Selection / active cell moves to O14 (as left by the author)
# ---------------------------------------------------------------------
$ws.Range("O14").Select()
